$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Program_sheet changes
# ---------------------------------------------------------------
$progSheet = $wb.Worksheets.Item("Program_sheet")

# Row 2: rename existing program entries
$progSheet.Range("A2").Value2 = "Yxyxyx"
$progSheet.Range("B2").Value2 = "ababab"

# New row 8: additional program entry (copy formatting from an existing
# data row first so the new row matches the sheet's row styling)
$progSheet.Range("A9:C9").Copy() | Out-Null
$progSheet.Range("A8:C8").PasteSpecial(-4122) | Out-Null
$progSheet.Range("A8").Value2 = "8playwright8"
$progSheet.Range("B8").Value2 = "auto"
$progSheet.Range("C8").Value2 = "active"

# Row 11: update status value
$progSheet.Range("B11").Value2 = "%G"

# New row 12: additional program entry (copy formatting from an existing
# data row first so the new row matches the sheet's row styling)
$progSheet.Range("A11:C11").Copy() | Out-Null
$progSheet.Range("A12:C12").PasteSpecial(-4122) | Out-Null
$progSheet.Range("A12").Value2 = "ML"
$progSheet.Range("B12").Value2 = "Machines"
$progSheet.Range("C12").Value2 = "Active"

$excel.CutCopyMode = 0

# Row 14: clear the long multi-line note in D14
$progSheet.Range("D14").ClearContents()

# ---------------------------------------------------------------
# Batch sheet changes
# ---------------------------------------------------------------
$batchSheet = $wb.Worksheets.Item("Batch")

# Rename shared batch name used across many rows (every cell that used to
# read "MobileSeleniumAuto" needs to be updated individually)
$batchSheet.Range("B2").Value2 = "activa"
$batchSheet.Range("B3").Value2 = "activa"
$batchSheet.Range("B4").Value2 = "activa"
$batchSheet.Range("B5").Value2 = "activa"
$batchSheet.Range("B7").Value2 = "activa"
$batchSheet.Range("B8").Value2 = "activa"
$batchSheet.Range("B9").Value2 = "activa"
$batchSheet.Range("B10").Value2 = "activa"

# Increase "Number of Classes" values (wait time for add new batch popup)
$batchSheet.Range("C2").Value2 = 13
$batchSheet.Range("C4").Value2 = 14
$batchSheet.Range("C5").Value2 = 15
